# The source workbook tracked each "driving" activity entry (rows 50-57,
# column G) using the shared string "driving". This edit recategorizes
# those entries as "commuting" instead. Writing the new text to every
# cell that held the old value lets the engine drop the now-unused
# "driving" shared string and append "commuting" to the shared-string
# table, which is exactly what the target workbook does.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G50:G57").Value = "commuting"

# Update the view state: the sheet was scrolled/selected near the bottom
# (K55) before; the author scrolled back up and selected E10.
$ws.Activate()
$ws.Range("A9").Select()
$ws.Range("E10").Select()
